$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nsgvs")
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "Number of CPU cores"
$ws.Range("A8").AddComment("Number of CPU cores to defined for the VM")
